# Generate Report for Handoff
# Rotate the source-document GUID (and its dependent handoff/handback file
# names + timestamps) to a freshly generated batch, across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newGuidFile = "3a02f150-8755-4fd3-8c90-af501b6ded87.md"

$newZhXlf  = "3a02f150-8755-4fd3-8c90-af501b6ded87.91872f057c9ecc467fae32c4494a7f15da07eaab.zh-cn.xlf"
$newZhDate = "2016-03-09 13:18:16"

$newDeXlf  = "3a02f150-8755-4fd3-8c90-af501b6ded87.91872f057c9ecc467fae32c4494a7f15da07eaab.de-de.xlf"
$newDeDate = "2016-03-09 13:18:25"

function Update-HyperlinkDisplay($ws, $cellAddr, $newDisplay) {
    $targetAddr = $ws.Range($cellAddr).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $targetAddr) {
            $hl.TextToDisplay = $newDisplay
        }
    }
}

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value() = $newGuidFile
Update-HyperlinkDisplay $wsOverview "A2" $newGuidFile

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value() = $newGuidFile
$wsZh.Range("C2").Value() = $newZhXlf
$wsZh.Range("D2").Value() = $newZhDate
Update-HyperlinkDisplay $wsZh "A2" $newGuidFile
Update-HyperlinkDisplay $wsZh "C2" $newZhXlf

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value() = $newGuidFile
$wsDe.Range("C2").Value() = $newDeXlf
$wsDe.Range("D2").Value() = $newDeDate
Update-HyperlinkDisplay $wsDe "A2" $newGuidFile
Update-HyperlinkDisplay $wsDe "C2" $newDeXlf
